$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, not auto-converted to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.741.20'
$ws.Range("E2").Value = '  +2.90%  '

$ws.Range("D3").Value = '2.010.96'
$ws.Range("E3").Value = '  +6.61%  '

$ws.Range("D4").Value = '1.015'
$ws.Range("E4").Value = '  +1.34%  '

$ws.Range("D5").Value = '330.82'
$ws.Range("E5").Value = '  +1.68%  '

$ws.Range("D6").Value = '1.009'
$ws.Range("E6").Value = '  +0.75%  '

$ws.Range("D7").Value = '0.4706'
$ws.Range("E7").Value = '  +2.48%  '

$ws.Range("D8").Value = '0.3972'
$ws.Range("E8").Value = '  +2.12%  '

$ws.Range("D9").Value = '47.19'
$ws.Range("E9").Value = '  +1.34%  '

$ws.Range("D10").Value = '0.08006'
$ws.Range("E10").Value = '  +1.35%  '

$ws.Range("D11").Value = '1.012'
$ws.Range("E11").Value = '  +2.67%  '

$ws.Range("D12").Value = '22.84'
$ws.Range("E12").Value = '  +4.74%  '

$ws.Range("D13").Value = '2.067.42'
$ws.Range("E13").Value = '  +8.56%  '

$ws.Range("D14").Value = '7.285'
$ws.Range("E14").Value = '  +3.78%  '

$ws.Range("D15").Value = '5.916'
$ws.Range("E15").Value = '  +4.08%  '

$ws.Range("D16").Value = '0.07196'
$ws.Range("E16").Value = '  +3.68%  '

$ws.Range("D17").Value = '89.53'
$ws.Range("E17").Value = '  +1.38%  '

$ws.Range("D18").Value = '1.009'
$ws.Range("E18").Value = '  +0.68%  '

$ws.Range("D19").Value = '0.00001004'
$ws.Range("E19").Value = '  +0.74%  '

$ws.Range("D20").Value = '17.53'
$ws.Range("E20").Value = '  +3.12%  '

$ws.Range("D21").Value = '1.007'
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").Value = '29.786.90'
$ws.Range("E22").Value = '  +3.06%  '

$ws.Range("D23").Value = '5.585'
$ws.Range("E23").Value = '  +5.75%  '

$ws.Range("D24").Value = '11.32'
$ws.Range("E24").Value = '  +3.24%  '

$ws.Range("D25").Value = '2.295.77'
$ws.Range("E25").Value = '  +7.90%  '

$ws.Range("D26").Value = '2.153'
$ws.Range("E26").Value = '  +3.30%  '

$ws.Range("D27").Value = '159.32'
$ws.Range("E27").Value = '  +2.75%  '

$ws.Range("D28").Value = '19.82'
$ws.Range("E28").Value = '  +2.61%  '

$ws.Range("D29").Value = '6.029'
$ws.Range("E29").Value = '  +0.75%  '

$ws.Range("D30").Value = '121.36'
$ws.Range("E30").Value = '  +3.27%  '

$ws.Range("D31").Value = '1.969'
$ws.Range("E31").Value = '  +2.00%  '

$ws.Range("D32").Value = '0.09499'
$ws.Range("E32").Value = '  +1.83%  '

$ws.Range("D33").Value = '0.9012'
$ws.Range("E33").Value = '  -0.28%  '

$ws.Range("D34").Value = '5.334'
$ws.Range("E34").Value = '  +0.96%  '

$ws.Range("D35").Value = '1.354'
$ws.Range("E35").Value = '  +2.02%  '

$ws.Range("D36").Value = '3.205'
$ws.Range("E36").Value = '  -1.81%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.05875'
$ws.Range("E37").Value = '  +2.02%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.185'
$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02142'
$ws.Range("E39").Value = '  +3.47%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.000003399'
$ws.Range("E40").Value = '  +103.49%  '

$ws.Range("D41").Value = '7.974'
$ws.Range("E41").Value = '  +4.40%  '

$ws.Range("D42").Value = '1.006'
$ws.Range("E42").Value = '  +0.48%  '

$ws.Range("D43").Value = '0.5808'
$ws.Range("E43").Value = '  +2.47%  '

$ws.Range("D44").Value = '0.1832'
$ws.Range("E44").Value = '  +3.80%  '

$ws.Range("D45").Value = '9.949'
$ws.Range("E45").Value = '  +2.68%  '

$ws.Range("D46").Value = '12.15'
$ws.Range("E46").Value = '  +2.47%  '

$ws.Range("D47").Value = '0.5424'
$ws.Range("E47").Value = '  +1.23%  '

$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '2.675'
$ws.Range("E48").Value = '  +4.93%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '2.176'
$ws.Range("E49").Value = '  -3.70%  '

$ws.Range("D50").Value = '0.07025'
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").Value = '1.886'
$ws.Range("E51").Value = '  +1.92%  '
